# Revert "adding term 2.0 now utf-8"
# Restores the workbook to the prior 1.1.0 term state:
#  - delete the extra "Include from FSIII 2" worksheet
#  - restore Version / Date / Contact metadata values
#  - restore the "descendent-of" value back to "D"

$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# --- Include from FSIII sheet updates ---------------------------------------
$inc = $wb.Worksheets.Item("Include from FSIII")
$inc.Range("C2").Value = "D"

# --- Remove the duplicated "Include from FSIII 2" sheet ---------------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Include from FSIII 2").Delete()
$excel.DisplayAlerts = $true
